$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A8: "N° Registro" -> "N° Pago", alignment center -> left (wrap stays on) ---
$ws.Range("A8").Value = "N° Pago"
$ws.Range("A8:D8").HorizontalAlignment = -4131   # xlLeft

# --- A22: "Monto Empeño:" -> "N° Registro" (style/formatting unchanged) ---
$ws.Range("A22").Value = "N° Registro"

# --- D22: set numeric value 0, accounting number format, bold (not underlined) font ---
$ws.Range("D22").Value = 0
$ws.Range("D22").NumberFormat = "_(* #,##0.00_);_(* \(#,##0.00\);_(* ""-""??_);_(@_)"
$ws.Range("D22").Font.Underline = -4142   # xlUnderlineStyleNone
$ws.Range("D22").Font.Bold = $true

# --- Selection moves from E7 to B16 ---
$ws.Range("B16").Select()

# --- Window scroll position ---
$excel.ActiveWindow.Top = 6264
